$wb = $excel.ActiveWorkbook

# --- "db" sheet: add a new "devPort" column between devIp and devUser ---
$wsDb = $wb.Worksheets.Item("db")
$wsDb.Columns("F").Insert()
$wsDb.Range("F2").Value = "devPort"

# --- selection / active-sheet bookkeeping (mirrors what Excel records on save) ---
# Set selections on the other two sheets first ...
$wsLib = $wb.Worksheets.Item("knihovny")
$wsLib.Range("B3").Select()

$wsRes = $wb.Worksheets.Item("zdroje")
$wsRes.Range("B6").Select()

# ... then finish on "db" so it ends up the active/selected tab, matching the diff.
$wsDb.Range("F2").Select()
